$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Column width changes - Overview sheet columns E and F (zh-cn, de-de)
$wsOverview.Range("E:E").ColumnWidth = 29.9777047293527
$wsOverview.Range("F:F").ColumnWidth = 29.9777047293527

# zh-cn sheet: column C (Status) and column P (Error Detail)
$wsZhCn.Range("C:C").ColumnWidth = 29.9777047293527
$wsZhCn.Range("P:P").ColumnWidth = 13.7470528738839

# de-de sheet: column C (Status) and column P (Error Detail)
$wsDeDe.Range("C:C").ColumnWidth = 29.9777047293527
$wsDeDe.Range("P:P").ColumnWidth = 13.7470528738839

# Update cell values - zh-cn sheet row 2
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-21 00:55:24"
$wsZhCn.Range("P2").Value = ""

# Update cell values - de-de sheet row 2
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-21 00:55:31"
$wsDeDe.Range("P2").Value = ""
